$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting "Maximum Gap" row down to row 9.
$ws.Rows.Item(8).Insert()

# Row 8 (new): duplicate Spiral Matrix entry logged against Bosscoder Academy
$ws.Range("C8").Value = "Spiral Matrix"
$ws.Range("D8").Value = "Bosscoder Academy"

# Rows 10-14: new problems solved in March
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "Maximum Subarray Sum"
$ws.Range("D10").Value = "Bosscoder Academy"

$ws.Range("C11").Value = "Maximum Product of Two Elements in an Array"
$ws.Range("D11").Value = "Bosscoder Academy"

$ws.Range("C12").Value = "Find N Unique Integers Sum up to Zero"
$ws.Range("D12").Value = "Bosscoder Academy"

$ws.Range("C13").Value = "Maximum Number of Pairs in Array"
$ws.Range("D13").Value = "Bosscoder Academy"

# Row 7: Spiral Matrix platform changes to IntelliJ IDEA (PC)
$ws.Range("D7").Value = "IntelliJ IDEA (PC)"

$ws.Range("C14").Value = "Counting Sort"
$ws.Range("D14").Value = "IntelliJ IDEA (PC)"

# Widen column C to fit the longer problem names
$ws.Columns.Item(3).ColumnWidth = 37.86

# Match the author's final selection
$ws.Range("H16").Select()
